# Orasi "AddEmployee" test fixture - roll the sample tester from
# "XIII" to "XIV" (data correction for the logout-enabled test run).
#
# Update the three related cells together (in one pass) so the shared
# string table keeps reusing the same <si> slots it already had for the
# old values, instead of retiring them and appending new ones at the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test.TestersonXIV"
$ws.Range("C2").Value = "TestersonXIV"
$ws.Range("M2").Value = "test.testersonxiv@orasi.com"

# The header row's selection (previously parked on G8 from a prior
# session) is cleared back to the sheet's default top-left cell.
[void]$ws.Range("A1").Select()

# Columns A, C and M re-measure their "best fit" width now that the
# text they hold has changed length.
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(3).ColumnWidth = 11.833333333333332
$ws.Columns.Item(13).ColumnWidth = 26.166666666666668
